$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at the very top; everything else shifts down by one.
$ws.Rows("1:1").Insert()

# Populate the new row 1 with the new transaction data.
$nbsp = [char]0x00A0
$concepto = "$nbsp$nbsp" + "TRANSFERENCIA INTERNET"
$monto = "10.00$nbsp$nbsp"

$ws.Range("A1").NumberFormat = "m/d/yy"
$ws.Range("A1").Value = (Get-Date -Year 2014 -Month 3 -Day 26).Date
$ws.Range("B1").Value = $concepto
$ws.Range("C1").Value = "C"
$ws.Range("D1").Value = "0004478289"
$ws.Range("E1").Value = "AG. NORTE"
$ws.Range("F1").Value = $monto
$ws.Range("G1").Value = "13.40"

# Re-create the shared CONCATENATE formula for row 1 (the insert pushed the
# previous H1 shared-formula owner down to H2, referencing row 2 now).
$f = $ws.Range("H2").Formula
$newF = $f -replace '([A-H])2(?!\d)', '${1}1'
$ws.Range("H1").Formula = $newF

# Restore the default selection to just H1 (matches the saved view state).
$ws.Range("H1").Select()
